# Apply "add prelim component for SDRS cola policy" edit:
#  - populate row 13 on paramlist with a new "cola_SDRS" policy row (clone of row 3
#    with cola_baseline -> cola_SDRS (A) and cola_type constant -> SDRS (E))
#  - flip B3 (include flag for cola_baseline) from TRUE to FALSE
#  - extend the H-column list data validation to cover the new row (H10:H12 -> H10:H13)
#  - make "paramlist" the active/selected sheet (was "Global_paramlist"), and update
#    the remembered selections on both sheets

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("paramlist")

# --- B3: include flag flips from TRUE to FALSE ---
$ws.Range("B3").Value = $false

# --- Row 13: new "cola_SDRS" policy (clone of row 3 with A/E swapped) ---
$ws.Range("A13").Value = "cola_SDRS"
$ws.Range("B13").Value = $true
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = "SDRS"
$ws.Range("F13").Value = "fixed"
$ws.Range("G13").Value = 0.015
$ws.Range("H13").Value = "ALpct"
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = 1
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.06
$ws.Range("M13").Value = 0.06
$ws.Range("R13").Value = 0.02
$ws.Range("S13").Value = 0
$ws.Range("T13").Value = 0.02
$ws.Range("U13").Value = 0
$ws.Range("V13").Value = 0.9
$ws.Range("W13").Value = 0.02
$ws.Range("X13").Value = 0
$ws.Range("Y13").Value = 0.9
$ws.Range("Z13").Value = 0.01
$ws.Range("AA13").Value = 0.001
$ws.Range("AB13").Value = 0.02
$ws.Range("AC13").Value = 0
$ws.Range("AD13").Value = 1
$ws.Range("AE13").Formula = "=0.15/20"
$ws.Range("AF13").Value = 0.001
$ws.Range("AG13").Value = 0.02
$ws.Range("AH13").Value = 0
$ws.Range("AI13").Value = 1
$ws.Range("AJ13").Value = 1
$ws.Range("AN13").Value = 0.06
$ws.Range("AO13").Value = 0.5
$ws.Range("AP13").Value = 0.5
$ws.Range("AQ13").Value = 0.04
$ws.Range("AR13").Value = 0.08
$ws.Range("AS13").Value = 0.04
$ws.Range("AT13").Value = 0.11
$ws.Range("AU13").Value = 0.04
$ws.Range("AV13").Value = 0.08
$ws.Range("AW13").Value = 0.7
$ws.Range("AX13").Value = 1
$ws.Range("AY13").Value = 15
$ws.Range("AZ13").Value = "closed"
$ws.Range("BA13").Value = "cd"
$ws.Range("BB13").Value = 0.02
$ws.Range("BC13").Value = "method1"
$ws.Range("BD13").Value = 5
$ws.Range("BE13").Value = 0.0822
$ws.Range("BF13").Value = 0.12
$ws.Range("BG13").Value = 0.075
$ws.Range("BH13").Value = 0.02
$ws.Range("BI13").Value = "constant"
$ws.Range("BJ13").Value = 0.02
$ws.Range("BK13").Value = 0.01
$ws.Range("BL13").Value = "AL_pct"
$ws.Range("BM13").Value = 0.75
$ws.Range("BN13").Value = $true
$ws.Range("BO13").Value = $true

# --- Extend the "H3 H5:H8 H10:H12" list validation to include H13 ---
# (Add across the contiguous block, then drop the two rows that must stay
# unvalidated, which the engine re-collapses back into one merged sqref.)
$ws.Range("H3").Validation.Delete()
$ws.Range("H5:H8").Validation.Delete()
$ws.Range("H10:H12").Validation.Delete()
$ws.Range("H3:H13").Validation.Add(3, 1, 1, "preSet, ALpct,MApct")
$ws.Range("H4").Validation.Delete()
$ws.Range("H9").Validation.Delete()

# --- Make "paramlist" the active sheet/tab, with B17 selected ---
$ws.Activate()
$ws.Range("B17").Select()

# --- Global_paramlist keeps its own remembered selection (K13) ---
$wsGlobal = $wb.Worksheets.Item("Global_paramlist")
$wsGlobal.Range("K13").Select()
